$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 83, pushing the existing rows 83:93 down to 84:94.
$ws.Rows.Item(83).Insert()

# Populate the newly inserted row 83 with the new weekly price observation
# (Ají / Inferno / Primera, Agrícola del Norte S.A. de Arica).
$ws.Cells.Item(83, 1).Value = 1
$ws.Cells.Item(83, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(83, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(83, 4).Value = 44785
$ws.Cells.Item(83, 5).Value = 15
$ws.Cells.Item(83, 6).Value = 100112021
$ws.Cells.Item(83, 7).Value = "Ají"
$ws.Cells.Item(83, 8).Value = "Inferno"
$ws.Cells.Item(83, 9).Value = "Primera"
$ws.Cells.Item(83, 10).Value = 130
$ws.Cells.Item(83, 11).Value = 11000
$ws.Cells.Item(83, 12).Value = 12000
$ws.Cells.Item(83, 13).Value = 11500
$ws.Cells.Item(83, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(83, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(83, 16).Value = 767
$ws.Cells.Item(83, 17).Value = 15
$ws.Cells.Item(83, 18).Value = "Hortaliza"
